# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# The metadata sheet re-classifies several columns that used to be tagged
# "iaest-dimension:*" as "iaest-measure:*" (they are now measures, not
# dimensions), re-points the "municipio-nombre" column at the SDMX
# reference-area dimension, gives it a dedicated "URI-Municipio" concept
# type (instead of the generic skos:Concept), and drops the leftover
# per-dimension mapping-file names from row 5 (only the ccaa-nombre mapping
# survives).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: iaest-dimension:... -> iaest-measure:... reclassifications,
#     plus municipio-nombre (L2) now pointing at sdmx-dimension:refArea.
$ws.Range("B2").Value = "iaest-measure:la-otra-actividad-del-titular-no-esta-relacionada-con-la-explotacion"
$ws.Range("D2").Value = "iaest-measure:el-titular-ha-realizado-actividades-complementarias-como-actividad-secundaria"
$ws.Range("L2").Value = "sdmx-dimension:refArea"
$ws.Range("N2").Value = "iaest-measure:menos-de-25-anos"
$ws.Range("U2").Value = "iaest-measure:la-otra-actividad-del-titular-esta-relacionada-con-la-explotacion"
$ws.Range("W2").Value = "iaest-measure:el-titular-no-ha-realizado-actividades-complementarias"
$ws.Range("AA2").Value = "iaest-measure:el-titular-ha-realizado-actividades-complementarias-como-actividad-principal"

# --- Row 3: "dim" -> "medida" for the same columns (L3 stays "dim" since
#     municipio-nombre is still a dimension, just a different one).
$ws.Range("B3").Value = "medida"
$ws.Range("D3").Value = "medida"
$ws.Range("L3").Value = "dim"
$ws.Range("N3").Value = "medida"
$ws.Range("U3").Value = "medida"
$ws.Range("W3").Value = "medida"
$ws.Range("AA3").Value = "medida"

# --- Row 4: "skos:Concept" -> "xsd:int" for the now-measure columns;
#     municipio-nombre (L4) gets its own concept type "URI-Municipio"
#     instead of the generic "skos:Concept".
$ws.Range("B4").Value = "xsd:int"
$ws.Range("D4").Value = "xsd:int"
$ws.Range("L4").Value = "URI-Municipio"
$ws.Range("N4").Value = "xsd:int"
$ws.Range("U4").Value = "xsd:int"
$ws.Range("W4").Value = "xsd:int"
$ws.Range("AA4").Value = "xsd:int"

# --- Row 5: only the ccaa-nombre mapping file reference remains; the other
#     per-dimension mapping file names are cleared out.
$ws.Range("B5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("U5").ClearContents()
$ws.Range("W5").ClearContents()
$ws.Range("AA5").ClearContents()
